$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optical_Power")

# Row 4
$ws.Range("A4").Value = "'6139"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "'6/12/2025"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "ASAMBLEA AV. 1440"
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 807458556
$ws.Range("F4").Value = "Optical Power"
$ws.Range("G4").Value = "Pendiente"
$ws.Range("H4").Value = "Cables en panza"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = '{"direccionesNormalizadas": [{"altura": 1440, "cod_calle": 1128, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.444800", "y": "-34.636534"}, "direccion": "ASAMBLEA AV. 1440, CABA", "nombre_calle": "ASAMBLEA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K4").Value = -58.4448
$ws.Range("L4").Value = -34.636534

# Row 5
$ws.Range("A5").Value = "'6135"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "'6/12/2025"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "JUJUY AV. 2179"
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 807458571
$ws.Range("F5").Value = "Optical Power"
$ws.Range("G5").Value = "Pendiente"
$ws.Range("H5").Value = "Cable en panza"
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = '{"direccionesNormalizadas": [{"altura": 2179, "cod_calle": 10013, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.399222", "y": "-34.635321"}, "direccion": "JUJUY AV. 2179, CABA", "nombre_calle": "JUJUY AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K5").Value = -58.399222
$ws.Range("L5").Value = -34.635321
